$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "27.434.34"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.563.28"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue "D5" "208.29"
$ws.Range("E5").Value = "  +1.42%  "
Set-TextValue "D6" "0.498"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  -0.12%  "
Set-TextValue "D8" "21.89"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  -0.76%  "
Set-TextValue "D10" "0.0589"
$ws.Range("E10").Value = "  +0.31%  "
Set-TextValue "D11" "0.0867"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.788.06"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "1.568.72"
$ws.Range("E13").Value = "  -0.24%  "
Set-TextValue "D14" "3.82"
$ws.Range("E14").Value = "  -0.30%  "
Set-TextValue "D15" "0.516"
$ws.Range("E15").Value = "  -2.19%  "
Set-TextValue "D16" "63.34"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "27.440.72"
$ws.Range("E17").Value = "  -0.12%  "
Set-TextValue "D18" "213.28"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  -0.09%  "
Set-TextValue "D20" "7.24"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -0.13%  "
Set-TextValue "D22" "4.11"
$ws.Range("E22").Value = "  -0.18%  "
Set-TextValue "D23" "9.53"
$ws.Range("E23").Value = "  +0.92%  "
Set-TextValue "D24" "2.02"
$ws.Range("E24").Value = "  +1.89%  "
Set-TextValue "D25" "153.00"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -0.12%  "
Set-TextValue "D27" "6.72"
$ws.Range("E27").Value = "  +0.64%  "
Set-TextValue "D28" "15.00"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  -1.42%  "
Set-TextValue "D30" "1.15"
$ws.Range("E30").Value = "  +0.06%  "
Set-TextValue "D31" "0.0469"
$ws.Range("E31").Value = "  +1.55%  "
Set-TextValue "D32" "3.19"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").Value = "1.359.08"
$ws.Range("E33").Value = "  +0.20%  "
Set-TextValue "D34" "2.94"
$ws.Range("E34").Value = "  +0.49%  "
Set-TextValue "D35" "1.53"
$ws.Range("E35").Value = "  +1.64%  "
Set-TextValue "D36" "0.973"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +2.24%  "
Set-TextValue "D39" "0.531"
$ws.Range("E39").Value = "  -0.59%  "
Set-TextValue "D40" "0.820"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("E41").Value = "  -0.16%  "
Set-TextValue "D42" "0.975"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  +1.75%  "
Set-TextValue "D44" "64.08"
$ws.Range("E44").Value = "  +1.23%  "
Set-TextValue "D45" "5.27"
$ws.Range("E45").Value = "  +0.80%  "
Set-TextValue "D46" "2.14"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "1.700.29"
$ws.Range("E47").Value = "  -0.33%  "
Set-TextValue "D48" "85.36"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -1.54%  "
Set-TextValue "D50" "0.0952"
$ws.Range("E50").Value = "  -1.26%  "
Set-TextValue "D51" "0.0494"
$ws.Range("E51").Value = "  -0.40%  "
